$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values after the cyclic rotation of rows 2-5 (D,I,J,K,L,M,N,P,Q columns)
$data = @{
    2 = @{ D = 44370; I = "Segunda"; J = 100; K = 1000; L = 1200; M = 1080; N = "`$/docena de matas"; P = 180;  Q = 6 }
    3 = @{ D = 44623; I = "Primera"; J = 300; K = 1800; L = 2000; M = 1900; N = "`$/paquete";          P = 1900; Q = 1 }
    4 = @{ D = 44377; I = "Segunda"; J = 550; K = 2000; L = 2800; M = 2364; N = "`$/docena de matas"; P = 394;  Q = 6 }
    5 = @{ D = 44267; I = "Primera"; J = 120; K = 1500; L = 1800; M = 1650; N = "`$/docena de matas"; P = 275;  Q = 6 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    $ws.Range("D$row").Value = $rowData.D
    $ws.Range("I$row").Value = $rowData.I
    $ws.Range("J$row").Value = $rowData.J
    $ws.Range("K$row").Value = $rowData.K
    $ws.Range("L$row").Value = $rowData.L
    $ws.Range("M$row").Value = $rowData.M
    $ws.Range("N$row").Value = $rowData.N
    $ws.Range("P$row").Value = $rowData.P
    $ws.Range("Q$row").Value = $rowData.Q
}
